# Applies the "Fixed update to excel issue" commit:
#   1. Renames the "Requested quantity" header on the "Weekly Quantity"
#      sheet to "Weekly_PO_Qty".
#   2. Renames the "Requested quantity" header on the "Monthly Trend"
#      sheet to "Monthly_PO_Qty".
#   3. Adds a new "PO Forecast" sheet (ds / PO_Forecast / yhat_lower /
#      yhat_upper) at the end of the workbook, populated with the
#      forecast series.

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$wsMonthly = $wb.Worksheets.Item(2)   # "Monthly Trend"

# --- 1 & 2: header renames -------------------------------------------------
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: new "PO Forecast" sheet --------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row values
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the header formatting used on the other sheets (bold, centered,
# bordered) by copying the existing header cell format rather than
# re-creating a new style.
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$poForecastData = @(
    @(44983.99999999999, 0, -9.328610526082427, 6.588067709779273),
    @(44990.99999999999, 0, -8.405298047685031, 6.139946359187591),
    @(45018.99999999999, 0, -7.103712203244884, 7.495214014083538),
    @(45039.99999999999, 1, -6.210786731604746, 9.256957889919248),
    @(45046.99999999999, 2, -6.043220605962373, 9.582152657977653),
    @(45053.99999999999, 2, -5.269266416571188, 9.552649668880127),
    @(45060.99999999999, 2, -4.723004962107746, 9.318089730255483),
    @(45088.99999999999, 4, -4.062784501879174, 11.35992379942413),
    @(45109.99999999999, 5, -3.193902976796284, 12.40727186514775),
    @(45116.99999999999, 5, -2.745794486172179, 12.66582088440877),
    @(45130.99999999999, 6, -1.373244624311145, 13.92853359343149),
    @(45137.99999999999, 6, -1.369596974554165, 13.88683037856266),
    @(45151.99999999999, 7, -0.5667650920930954, 14.1556563131828),
    @(45165.99999999999, 7, -0.3366734252903113, 14.8905058682286),
    @(45179.99999999999, 8, 0.1249815437621832, 15.59379818703069),
    @(45186.99999999999, 8, 0.7409268646165125, 15.47597811415699),
    @(45193.99999999999, 9, 1.233994461064829, 16.3063332581716),
    @(45200.99999999999, 9, 1.409886215454121, 16.87156784918817),
    @(45207.99999999999, 10, 1.804293909655246, 17.18855183908811),
    @(45214.99999999999, 10, 2.013256891216585, 17.98638676558129),
    @(45221.99999999999, 10, 2.491392306672088, 17.86931821250781),
    @(45256.99999999999, 12, 3.976514962469605, 19.60629096411755),
    @(45312.99999999999, 15, 6.666408471373924, 21.59294214326417),
    @(45340.99999999999, 16, 8.650721604159374, 23.63115279117789),
    @(45354.99999999999, 17, 9.544514022518749, 24.63905621786057),
    @(45361.99999999999, 17, 9.258362548266311, 24.57768394099404),
    @(45368.99999999999, 17, 9.298648075449613, 24.82719789139326),
    @(45375.99999999999, 18, 10.61401337999707, 25.48158984623494),
    @(45382.99999999999, 18, 10.34991847544637, 25.30276099256608),
    @(45389.99999999999, 18, 10.53684903816501, 25.7534097386284),
    @(45396.99999999999, 19, 10.38622753865239, 25.9543418011892),
    @(45403.99999999999, 19, 11.00977308022365, 26.4393451885655),
    @(45410.99999999999, 19, 11.690369109223, 26.88986555005556)
)

$r = 2
foreach ($row in $poForecastData) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Match the date-column formatting used in column A of the other sheets.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A34").PasteSpecial(-4122)  # xlPasteFormats

[void]$wsForecast.Range("A1").Select()
